$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (shifts existing columns A-D to B-E)
$ws.Columns.Item(1).Insert()

# Row 1 headers
$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# Row 2 values
$ws.Range("A2").Value = "CasesTab"

$ws.Range("B2").Value = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = `"WHITE`"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS ``Case ID``,
    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,
    COALESCE(a.arm_id, '') AS ``Arm``,
    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,
    COALESCE(c.disease, '') AS ``Diagnosis``,
    COALESCE(c.gender, '') AS ``Gender``,
    COALESCE(c.race, '') AS ``Race``,
    COALESCE(c.ethnicity, '') AS ``Ethnicity``"

$ws.Range("C2").Value = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = `"WHITE`"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"

# D2/E2 remain the Neo4jData/WebData filenames (already shifted by the column insert)
# Columns B-E keep their original widths automatically via the column insert/shift above.

# New column A width (auto-fit to the short "TabName"/"CasesTab" header text)
$ws.Columns.Item(1).ColumnWidth = 8

# Row height for row 2
$ws.Rows.Item(2).RowHeight = 174

# Apply wrap text style to B2 and C2 (style index 1 = wrapText alignment)
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true

# Selection
$ws.Range("B10").Select() | Out-Null
